$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as plain text (e.g. grouped
# thousands like "28.397.37", or numbers with significant trailing
# zeros like "53.60"). A bare .Value assignment lets Excel
# auto-convert numeric-looking text into a real number, which would
# silently corrupt the value (lose trailing zeros / change type).
# Forcing NumberFormat to Text ("@") before the write keeps the
# text exactly as scraped, then resetting the Style back to Normal
# avoids leaving a stray formatting change behind.

$ws.Range("D2").Value = '28.397.37'
$ws.Range("E2").Value = '  -2.46%  '
$ws.Range("D3").Value = '1.947.82'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4797'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4089'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08507'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.056'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.02%  '
$ws.Range("D13").Value = '2.000.06'
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.562'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.151'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.014'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001073'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06639'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.88%  '
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.836'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = '28.437.99'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.299'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").Value = '2.185.09'
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.176'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.834'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '124.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9839'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09661'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.697'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.620'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.435'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.166'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02331'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06184'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.248'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6220'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.011'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1919'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("E45").Value = '  +3.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5951'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.060'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.407'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06814'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000306'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.88%  '
